$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Fishing
$ws.Range("C2").Value = 0.7404891432017738
$ws.Range("D2").Value = 0.7404891432017738
$ws.Range("E2").Value = 3.6135431657076023
$ws.Range("F2").Value = 0.019710524742908027
$ws.Range("G2").Value = 0.0014

# Row 3 - Bioregion
$ws.Range("C3").Value = 6.836531002670082
$ws.Range("D3").Value = 2.278843667556694
$ws.Range("E3").Value = 11.12062214040009
$ws.Range("F3").Value = 0.18197648773233718

# Row 4 - Fishing:Bioregion
$ws.Range("C4").Value = 1.3023176368950224
$ws.Range("D4").Value = 0.43410587896500746
$ws.Range("E4").Value = 2.118410980807661
$ws.Range("F4").Value = 0.03466541574688594
$ws.Range("G4").Value = 0.0031

# Row 5 - Residuals
$ws.Range("C5").Value = 28.688872747407192
$ws.Range("D5").Value = 0.20492051962433708
$ws.Range("F5").Value = 0.7636475717778689

# Row 6 - Total
$ws.Range("C6").Value = 37.56821053017407
